# riv_data.xlsx revision:
# "revised obs. ref, river slope, extended cal period, fixed rbot"
#
# Data edits on the "model" sheet:
#  - E2 (observation reference level) revised from 20.6 to 0
#  - D3, D5, D8 (river slope values) revised to 0.06 and highlighted in red
#    (new font/style) to flag the fixed values
#  - dependent E-column formulas recompute automatically
#
# Plus cosmetic chart repositioning on the "2024", "2021" and "Comparaison "
# sheets, and an updated active-cell selection on "model".

$wb = $excel.ActiveWorkbook

# ---- sheet "model" ------------------------------------------------------
$ws = $wb.Worksheets.Item("model")

# revised observation reference level
$ws.Range("E2").Value = 0

# revised river slope values (D3, D5, D8) -> 0.06, flagged in red
$ws.Range("D3").Value = 0.06
$ws.Range("D3").Font.Color = 255

$ws.Range("D5").Value = 0.06
$ws.Range("D5").Font.Color = 255

$ws.Range("D8").Value = 0.06
$ws.Range("D8").Font.Color = 255

# updated active cell selection
[void]$ws.Range("D5").Select()

# ---- chart/group repositioning -------------------------------------------
$ws2024 = $wb.Worksheets.Item("2024")
$grp = $ws2024.Shapes.Item("Groupe 1")
$grp.Left = 535.55
$grp.Top = 278.82
$grp.Width = 520.83
$grp.Height = 324.56

$ws2021 = $wb.Worksheets.Item("2021")
$grp = $ws2021.Shapes.Item("Groupe 1")
$grp.Left = 855.3
$grp.Top = 8.7
$grp.Width = 462
$grp.Height = 312

$wsComp = $wb.Worksheets.Item("Comparaison ")
$grp1 = $wsComp.Shapes.Item("Groupe 1")
$grp1.Left = 431.63
$grp1.Top = 32.32
$grp1.Width = 524.37
$grp1.Height = 298.37

$grp2 = $wsComp.Shapes.Item("Groupe 2")
$grp2.Left = 9.17
$grp2.Top = 31.99
$grp2.Width = 459.52
$grp2.Height = 299.08
